# Apply the "add 2022-Q4 data" edit:
#  1. Insert a new worksheet "2022-Q4" right after "总计" holding the
#     fund-holder breakdown for the new quarter.
#  2. Prepend a matching summary row to the "总计" sheet (shifting the
#     previously-existing rows down by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Shift the existing "总计" summary rows down and insert the new one.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Capture the current (pre-edit) rows 2..5 before overwriting anything.
$rows = @()
for ($r = 2; $r -le 5; $r++) {
    $rows += , @(
        $total.Cells.Item($r, 2).Value2,
        $total.Cells.Item($r, 3).Value2,
        $total.Cells.Item($r, 4).Value2
    )
}

# Write the new first data row (2022-Q4) into row 2.
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 3
$total.Cells.Item(2, 4).Value = 0.09

# Extend the formatted range down to row 6 (copy row 5's styling, which
# already carries the correct bordered/centered style for column A).
$total.Rows.Item(5).Copy($total.Rows.Item(6))

# Push the previously-captured rows down one position (row 2->3, 3->4, ...).
for ($i = 0; $i -lt $rows.Length; $i++) {
    $destRow = $i + 3
    $total.Cells.Item($destRow, 2).Value = $rows[$i][0]
    $total.Cells.Item($destRow, 3).Value = $rows[$i][1]
    $total.Cells.Item($destRow, 4).Value = $rows[$i][2]
}

# Column A is just the zero-based row index; extend it to the new last row.
$total.Cells.Item(6, 1).Value = 4

# ---------------------------------------------------------------------
# 2) Add the new "2022-Q4" worksheet right after "总计". Cloning an
#    existing quarter sheet (instead of Worksheets.Add) means the new
#    sheet starts with identical header/row styling.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy($null, $total)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Extend formatting down to rows 3 and 4 (the template only had one data row).
$q4.Rows.Item(2).Copy($q4.Rows.Item(3))
$q4.Rows.Item(2).Copy($q4.Rows.Item(4))

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q4.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$data = @(
    @("970007", "华安证券汇赢增利一年持有混合B", "11.05", "22.39", "0.48", "0.0530", 6),
    @("970008", "华安证券汇赢增利一年持有混合C", "8.56", "22.39", "0.48", "0.0411", 6),
    @("970006", "华安证券汇赢增利一年持有混合A", "0.18", "22.39", "0.48", "0.0009", 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $q4.Cells.Item($r, 1).Value = $i

    # Columns B, D, E, F, G hold numeric-looking text ("970007", "11.05", ...)
    # that must stay plain text (inlineStr), matching the source data. Force
    # text via a temporary "@" number format, then ClearFormats() so the
    # cell falls back to the default (unstyled) xf instead of keeping a
    # bespoke "text" style index.
    $q4.Cells.Item($r, 2).NumberFormat = "@"
    $q4.Cells.Item($r, 2).Value = $row[0]
    $q4.Cells.Item($r, 2).ClearFormats()

    $q4.Cells.Item($r, 3).Value = $row[1]

    $q4.Cells.Item($r, 4).NumberFormat = "@"
    $q4.Cells.Item($r, 4).Value = $row[2]
    $q4.Cells.Item($r, 4).ClearFormats()

    $q4.Cells.Item($r, 5).NumberFormat = "@"
    $q4.Cells.Item($r, 5).Value = $row[3]
    $q4.Cells.Item($r, 5).ClearFormats()

    $q4.Cells.Item($r, 6).NumberFormat = "@"
    $q4.Cells.Item($r, 6).Value = $row[4]
    $q4.Cells.Item($r, 6).ClearFormats()

    $q4.Cells.Item($r, 7).NumberFormat = "@"
    $q4.Cells.Item($r, 7).Value = $row[5]
    $q4.Cells.Item($r, 7).ClearFormats()

    $q4.Cells.Item($r, 8).Value = $row[6]
}
